# Prefix each worksheet's own name (plus a space) onto the Step/command
# names found in column A of that worksheet (every data row below the
# "Name" header row, i.e. starting at row 2), for every "protocol" sheet
# in the workbook. The first five overview/reference sheets
# (MaddisonJourney, NRWaves, PersonalMaddison, PositiveSpin, ReEngagement)
# are left untouched.

$wb = $excel.ActiveWorkbook

$skipSheets = @(
    "MaddisonJourney",
    "NRWaves",
    "PersonalMaddison",
    "PositiveSpin",
    "ReEngagement"
)

foreach ($ws in $wb.Worksheets) {
    if ($skipSheets -contains $ws.Name) {
        continue
    }

    $prefix = $ws.Name

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value2
        if ($null -ne $val -and $val -ne "") {
            $text = [string]$val
            if (-not $text.StartsWith("$prefix ")) {
                $cell.Value2 = "$prefix $text"
            }
        }
    }
}
